# Fix formatting issues introduced when scraping floating point numbers
# and some provider names in the contratos sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Provider/company name fixes: stray commas used as separators were
#    scraped incorrectly and must be periods; "S.H." should read "SH".
$nameFixes = @{
    "E51"  = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
    "E77"  = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
    "E185" = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
    "E78"  = "BOFFELLI. MARIA INES"
    "E121" = "GIMENEZ ANIBAL. FALISTOCCO MARISA DANIELA SH"
    "E183" = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
    "E194" = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
}

foreach ($addr in $nameFixes.Keys) {
    $ws.Range($addr).Value = $nameFixes[$addr]
}

# 2) "Importe" column (H) amounts were scraped using the Argentine
#    locale (e.g. "9.120,00") instead of a plain decimal point format
#    (e.g. "9120.00"). Strip the thousands separators ('.') and turn the
#    decimal comma (',') into a decimal point ('.'). The values are kept
#    as text (not converted to real numbers), so force a text format
#    before writing them back.
$importeRange = $ws.Range("H2:H286")
$importeRange.NumberFormat = "@"

for ($r = 2; $r -le 286; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $old = $cell.Value2
    $new = $old.Replace(".", "").Replace(",", ".")
    $cell.Value = $new
}
